$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target table (18 data rows) -- a new bucket "UA5" (4 parcels) is
# inserted at the top of the table, the existing buckets (UA3, UA4, UA1,
# gUA2) shift down and their ranking values are recomputed; also the
# holding_id of the parcel with parcel_id 13 (UA1 bucket, ranking 7)
# changes from "A" to "B" because that holding now contains a
# high-ranking parcel that quickly fills the bucket.
$rows = @(
    @(1,  "UA5",  75,   "A", "UA5",  2),
    @(2,  "UA5",  7723, "B", "UA5",  3),
    @(3,  "UA5",  3,    "B", "UA5",  11),
    @(4,  "UA5",  988,  "C", "UA5",  10),
    @(5,  "UA3",  601,  "A", "UA3",  9),
    @(6,  "UA3",  2195, "A", "UA3",  12),
    @(7,  "UA3",  7723, "B", "UA3",  3),
    @(8,  "UA4",  3289, "B", "UA4",  8),
    @(9,  "UA4",  3,    "B", "UA4",  11),
    @(10, "UA4",  988,  "C", "UA4",  10),
    @(11, "UA1",  6235, "A", "UA1",  1),
    @(12, "UA1",  75,   "A", "UA1",  2),
    @(13, "UA1",  753,  "B", "UA1",  7),
    @(14, "UA1",  7723, "B", "UA1",  3),
    @(15, "gUA2", 75,   "A", "gUA2", 2),
    @(16, "gUA2", 2195, "A", "gUA2", 12),
    @(17, "gUA2", 1194, "B", "gUA2", 6),
    @(18, "gUA2", 3289, "B", "gUA2", 8)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
}

# Apply the same format (bold, bordered, centered/top-aligned) used by the
# existing column-A cells to the newly added column-A cells (rows 16-19).
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A16:A19").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
